$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "objectid" row (old row 2). This shifts id_punto..latitud up by one row.
$ws.Rows.Item(2).Delete()

# Update the "Descripción" (column B) and "Tipo" (column C) values to reflect the
# contextualized wording ("de la tarjeta Cívica" instead of "de Civica") and the
# corrected data types.

# Row 2: id_punto
$ws.Cells.Item(2, 2).Value = "Identificador del punto de recarga de la tarjeta Cívica"
$ws.Cells.Item(2, 3).Value = "Número"

# Row 3: nombreestablecimiento
$ws.Cells.Item(3, 2).Value = "Nombre del establecimiento disponible para recarga de la tarjeta Cívica"
$ws.Cells.Item(3, 3).Value = "Texto "

# Row 4: direccion
$ws.Cells.Item(4, 2).Value = "Dirección del establecimiento de recarga de la tarjeta Cívica"
$ws.Cells.Item(4, 3).Value = "Texto "

# Row 5: municipio
$ws.Cells.Item(5, 2).Value = "Municipio del establecimiento de recarga de la tarjeta Cívica"
$ws.Cells.Item(5, 3).Value = "Texto "

# Row 6: barrio
$ws.Cells.Item(6, 2).Value = "Barrio del establecimiento de recarga de la tarjeta Cívica"
$ws.Cells.Item(6, 3).Value = "Texto "

# Row 7: estado
$ws.Cells.Item(7, 2).Value = "Estado del establecimiento de recarga de la tarjeta Cívica"
$ws.Cells.Item(7, 3).Value = "Texto "

# Row 8: longitud
$ws.Cells.Item(8, 2).Value = "Longuitud de las localizaciones de los establecimientos de recarga de la tarjeta Cívica"
$ws.Cells.Item(8, 3).Value = "Número"

# Row 9: latitud
$ws.Cells.Item(9, 2).Value = "Latitud de las localizaciones de los establecimientos de recarga de la tarjeta Cívica"
$ws.Cells.Item(9, 3).Value = "Número"

# Minor cosmetic touch-ups matching the re-saved workbook: column widths and the
# active selection.
$ws.Columns.Item(1).ColumnWidth = 21.7
$ws.Columns.Item(2).ColumnWidth = 47.0

$ws.Range("B2").Select() | Out-Null
